# Agregado de menu (sin funcionalidad). Actualización de presentación.
#
# Reposiciona algunas formas del slide de portada (slide 1) y fusiona dos
# runs de texto en el slide 4 ("Tipos de Clientes") en un unico run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 (portada): ajustar la posicion horizontal de varias formas.
# Los valores se expresan en puntos (1 punto = 12700 EMU), que es la
# unidad que usa el modelo de objetos de PowerPoint para Shape.Left/Top.
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

# "1 Título" (marcador de título, placeholder ctrTitle) -> x: 611560 -> 687760 EMU
$slide1.Shapes.Item(1).Left = 687760 / 12700

# "2 Subtítulo" (marcador de subtitulo) -> x: 899592 -> 902668 EMU
$slide1.Shapes.Item(2).Left = 902668 / 12700

# "Picture 3" (imagen central) -> x: 2555776 -> 2833989 EMU
$slide1.Shapes.Item(3).Left = 2833989 / 12700

# "2 Subtítulo" (cuadro de texto "Grupo 5...") -> x: 827584 -> 866664 EMU
$slide1.Shapes.Item(4).Left = 866664 / 12700

# "1 Título" (cuadro de texto "Panic Dial Button") -> x: 763960 -> 687760 EMU
$slide1.Shapes.Item(5).Left = 687760 / 12700

# ---------------------------------------------------------------------
# Slide 4 ("Tipos de Clientes"): en el parrafo "Planes según Necesidades
# Específicas: Empresarial o Familiar" los runs " " y "Empresarial o
# Familiar" se combinan en un unico run " Empresarial o Familiar".
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$textBox = $slide4.Shapes.Item(6)
$tr = $textBox.TextFrame.TextRange

$fullText = $tr.Text
$target = " Empresarial o Familiar"
$idx0 = $fullText.IndexOf($target)
if ($idx0 -ge 0) {
    $run = $tr.Characters($idx0 + 1, $target.Length)
    $run.Text = $target
}
